# Apply the paragraph-formatting defaults that the "Normal" style picked
# up (left-to-right reading order, zero before/after spacing, and
# explicit left justification) in word/styles.xml.
#
# These three settings surface in the OOXML as, respectively:
#   <w:bidi w:val="0"/>
#   <w:spacing w:before="0" w:after="0"/>
#   <w:jc w:val="left"/>
# appended (in that order) to the Normal style's <w:pPr>, right after the
# pre-existing <w:widowControl/>.

$d = $word.ActiveDocument

$normal = $d.Styles("Normal")
$pf = $normal.ParagraphFormat

$pf.ReadingOrder = 0   # wdReadingOrderLtr -> <w:bidi w:val="0"/>
$pf.SpaceBefore = 0    # -> <w:spacing w:before="0" .../>
$pf.SpaceAfter = 0     # -> <w:spacing ... w:after="0"/>
$pf.Alignment = 0      # wdAlignParagraphLeft -> <w:jc w:val="left"/>
